$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.916.70"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.635.02"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.20"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5069"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2576"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06353"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07748"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.257"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.638.20"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5512"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "0.0₅7693"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.91"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "25.929.82"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.442"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.24"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.898"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.053"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.907"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.25"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1235"
$ws.Range("E26").Value = "  +5.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.810"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.54"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.244"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04883"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.251"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.188"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.542"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.377"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9047"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.566"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5493"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "1.122.15"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01556"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.570"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8057"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.50"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("D45").Value = "1.773.37"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4457"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.84"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9949"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05146"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.557"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("E51").Value = "  -0.26%  "
